$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the row labels to shortened forms ("S1" -> "S", "B1" -> "B")
$ws.Range("A5").Value = "B"
$ws.Range("A4").Value = "S"

# Update the active selection on the sheet
$ws.Range("A5").Select()
